$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047371756"
$ws.Range("D16").Value = "MIRLEY MONTALVO PEREZ"
$ws.Range("E16").Value = "2103"
$ws.Range("F16").Value = 10820
$ws.Range("G16").Value = 737717

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047371756"
$ws.Range("D17").Value = "MIRLEY MONTALVO PEREZ"
$ws.Range("E17").Value = "2102"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 737717

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047371756"
$ws.Range("D18").Value = "MIRLEY MONTALVO PEREZ"
$ws.Range("E18").Value = "2101"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 737717

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047371756"
$ws.Range("D19").Value = "MIRLEY MONTALVO PEREZ"
$ws.Range("E19").Value = "2012"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 828116

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047371756"
$ws.Range("D20").Value = "MIRLEY MONTALVO PEREZ"
$ws.Range("E20").Value = "2011"
$ws.Range("F20").Value = 33125
$ws.Range("G20").Value = 828116

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1047371756"
$ws.Range("D21").Value = "MIRLEY MONTALVO PEREZ"
$ws.Range("E21").Value = "2010"
$ws.Range("F21").Value = 33125
$ws.Range("G21").Value = 828116

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047371756"
$ws.Range("D22").Value = "MIRLEY MONTALVO PEREZ"
$ws.Range("E22").Value = "2009"
$ws.Range("F22").Value = 33125
$ws.Range("G22").Value = 828116

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1047371756"
$ws.Range("D23").Value = "MIRLEY MONTALVO PEREZ"
$ws.Range("E23").Value = "2008"
$ws.Range("F23").Value = 33125
$ws.Range("G23").Value = 828116

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1047371756"
$ws.Range("D24").Value = "MIRLEY MONTALVO PEREZ"
$ws.Range("E24").Value = "2007"
$ws.Range("F24").Value = 33125
$ws.Range("G24").Value = 828116

$ws.Range("B25").Value = "CE"
$ws.Range("C25").Value = "18140237"
$ws.Range("D25").Value = "CARMEN ELENA HERRERA GARCIA"
$ws.Range("E25").Value = "1801"
$ws.Range("F25").Value = 29509
$ws.Range("G25").Value = 737717

$ws.Range("B26").Value = "CE"
$ws.Range("C26").Value = "18140237"
$ws.Range("D26").Value = "CARMEN ELENA HERRERA GARCIA"
$ws.Range("E26").Value = "1712"
$ws.Range("F26").Value = 29509
$ws.Range("G26").Value = 737717

$ws.Range("B27").Value = "CE"
$ws.Range("C27").Value = "18140237"
$ws.Range("D27").Value = "CARMEN ELENA HERRERA GARCIA"
$ws.Range("E27").Value = "1711"
$ws.Range("F27").Value = 10820
$ws.Range("G27").Value = 737717

